$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: fill in Start/End times and add the "Moved graph view..." log entry ---
$ws.Range("B6").Value = 0.42722222222222223
$ws.Range("C6").Value = 0.49140046296296297
$ws.Range("F6").Value = "* Moved graph view into an intermediary visual element that holds everything (toolbar, blackboard etc)`n* Fixed blackboard visibility issue"

# --- Row 17: fill in Start/End times and add the "Added Properties..." log entry ---
$ws.Range("B17").Value = 0.55943287037037037
$ws.Range("C17").Value = 0.81623842592592588
$ws.Range("F17").Value = "* Added Properties (the checks & triggers I was talking about, and Actors which are basically a property that  holds which entity owns the dialogue line)`n* Fixed blackboard: you can now properly delete properties.`n* Added property (de-)serialization`n* Added node creation from property (either using the search window provider or by dragging a property from the blackboard)`n* Did some more research into GraphView and UIElements.`n* Experimented a bit with writing USS (CSS/UI Styling for UIElements)"
$ws.Range("F17").WrapText = $true

# Setting wrapped multi-line text can trigger row auto-fit; restore the
# original fixed row height (15) to match the source formatting.
$ws.Rows("6:6").RowHeight = 15
$ws.Rows("17:17").RowHeight = 15

# --- Hide the rows for days with no work logged (rows 3 and 7-16) ---
$ws.Rows("3:3").Hidden = $true
$ws.Rows("7:16").Hidden = $true

# --- Apply the duration number format ([h]:mm:ss) to the Total Time cell ---
$ws.Range("G2").NumberFormat = "[h]:mm:ss"

# --- Update the active selection to F17 ---
$ws.Range("F17").Select() | Out-Null
